$wb = $excel.ActiveWorkbook

# --- training_site_summary sheet ---
$ws1 = $wb.Worksheets.Item("training_site_summary")
$ws1.Range("B7").Value = 20
$ws1.Range("C7").Value = 690
$ws1.Range("B8").Value = 175
$ws1.Range("C8").Value = 4828

# --- training_subject_summary sheet ---
$ws2 = $wb.Worksheets.Item("training_subject_summary")

# Row 68
$ws2.Range("C68").Value = "In progress"
$ws2.Range("D68").Value = 4
$ws2.Range("E68").Value = 2
$ws2.Range("F68").Value = 0
$ws2.Range("G68").Value = 0
$ws2.Range("H68").Value = 4
$ws2.Range("I68").Value = "22-08-2023"
$ws2.Range("J68").Value = "not match"
$ws2.Range("U68").Value = 1

# Row 69
$ws2.Range("C69").Value = "In progress"
$ws2.Range("D69").Value = 8
$ws2.Range("E69").Value = 0
$ws2.Range("F69").Value = 0
$ws2.Range("G69").Value = 0
$ws2.Range("H69").Value = 8
$ws2.Range("I69").Value = "22-08-2023"
$ws2.Range("J69").Value = "30-08-2023"
$ws2.Range("K69").Value = "not match"
$ws2.Range("U69").Value = 2

# Row 70
$ws2.Range("C70").Value = "In progress"
$ws2.Range("D70").Value = 4
$ws2.Range("E70").Value = 0
$ws2.Range("F70").Value = 0
$ws2.Range("G70").Value = 0
$ws2.Range("H70").Value = 4
$ws2.Range("I70").Value = "22-08-2023"
$ws2.Range("U70").Value = 1

# Row 71
$ws2.Range("C71").Value = "In progress"
$ws2.Range("D71").Value = 2
$ws2.Range("E71").Value = 0
$ws2.Range("F71").Value = 0
$ws2.Range("G71").Value = 0
$ws2.Range("H71").Value = 2
$ws2.Range("I71").Value = "22-08-2023"
$ws2.Range("J71").Value = "not match"
$ws2.Range("K71").Value = "not match"
$ws2.Range("U71").Value = 1

# Row 72
$ws2.Range("C72").Value = "In progress"
$ws2.Range("D72").Value = 4
$ws2.Range("E72").Value = 4
$ws2.Range("F72").Value = 0
$ws2.Range("G72").Value = 0
$ws2.Range("H72").Value = 4
$ws2.Range("I72").Value = "22-08-2023"
$ws2.Range("U72").Value = 1

# Row 144
$ws2.Range("N144").Value = "not match"
$ws2.Range("O144").Value = "not match"
$ws2.Range("P144").Value = "not match"

# Row 145
$ws2.Range("L145").Value = "not match"
$ws2.Range("M145").Value = "not match"
$ws2.Range("N145").Value = "not match"

# Row 147
$ws2.Range("L147").Value = "not match"
$ws2.Range("M147").Value = "not match"
$ws2.Range("N147").Value = "not match"

# Row 148
$ws2.Range("J148").Value = "not match"
$ws2.Range("K148").Value = "not match"
$ws2.Range("L148").Value = "not match"

# Row 149
$ws2.Range("K149").Value = "not match"
$ws2.Range("L149").Value = "not match"
$ws2.Range("M149").Value = "not match"

# Row 150
$ws2.Range("J150").Value = "not match"

# Row 160
$ws2.Range("B160").Value = 68
$ws2.Range("D160").Value = 68
$ws2.Range("H160").Value = 10
$ws2.Range("S160").Value = "03-08-2023"
$ws2.Range("T160").Value = "15-08-2023"
$ws2.Range("U160").Value = 12

# Row 163
$ws2.Range("B163").Value = 76
$ws2.Range("D163").Value = 71
$ws2.Range("E163").Value = 5
$ws2.Range("H163").Value = 40
$ws2.Range("N163").Value = "01-08-2023"
$ws2.Range("O163").Value = "08-08-2023"
$ws2.Range("P163").Value = "15-08-2023"
$ws2.Range("Q163").Value = "22-08-2023"
$ws2.Range("R163").Value = "29-08-2023"
$ws2.Range("S163").Value = "05-09-2023"
$ws2.Range("U163").Value = 11

# Row 165
$ws2.Range("B165").Value = 28
$ws2.Range("D165").Value = 28
$ws2.Range("H165").Value = 14
$ws2.Range("M165").Value = "27-07-2023"
$ws2.Range("N165").Value = "10-08-2023"
$ws2.Range("U165").Value = 5

# Row 167
$ws2.Range("B167").Value = 51
$ws2.Range("D167").Value = 51
$ws2.Range("H167").Value = 9
$ws2.Range("P167").Value = "30-08-2023"
$ws2.Range("Q167").Value = "05-09-2023"
$ws2.Range("U167").Value = 9

# Row 168
$ws2.Range("B168").Value = 44
$ws2.Range("D168").Value = 44
$ws2.Range("H168").Value = 14
$ws2.Range("N168").Value = "27-07-2023"
$ws2.Range("O168").Value = "03-08-2023"
$ws2.Range("P168").Value = "14-08-2023"
$ws2.Range("U168").Value = 7

# Row 169
$ws2.Range("B169").Value = 36
$ws2.Range("D169").Value = 36
$ws2.Range("H169").Value = 18
$ws2.Range("M169").Value = "31-07-2023"
$ws2.Range("N169").Value = "08-08-2023"
$ws2.Range("O169").Value = "21-08-2023"
$ws2.Range("P169").Value = "04-09-2023"
$ws2.Range("U169").Value = 8

# Row 170
$ws2.Range("B170").Value = 27
$ws2.Range("D170").Value = 27
$ws2.Range("H170").Value = 20
$ws2.Range("J170").Value = "26-07-2023"
$ws2.Range("K170").Value = "02-08-2023"
$ws2.Range("L170").Value = "16-08-2023"
$ws2.Range("M170").Value = "28-08-2023"
$ws2.Range("U170").Value = 5

# Row 171
$ws2.Range("B171").Value = 22
$ws2.Range("D171").Value = 22
$ws2.Range("H171").Value = 17
$ws2.Range("J171").Value = "27-07-2023"
$ws2.Range("K171").Value = "10-08-2023"
$ws2.Range("L171").Value = "18-08-2023"
$ws2.Range("M171").Value = "04-09-2023"
$ws2.Range("U171").Value = 5

# Row 172
$ws2.Range("A172").Value = "292-016"
$ws2.Range("B172").Value = 14
$ws2.Range("C172").Value = "In progress"
$ws2.Range("D172").Value = 14
$ws2.Range("E172").Value = 0
$ws2.Range("F172").Value = 0
$ws2.Range("G172").Value = 0
$ws2.Range("H172").Value = 14
$ws2.Range("I172").Value = "03-08-2023"
$ws2.Range("J172").Value = "09-08-2023"
$ws2.Range("U172").Value = 2

# Row 173
$ws2.Range("A173").Value = "292-017"
$ws2.Range("B173").Value = 7
$ws2.Range("C173").Value = "In progress"
$ws2.Range("I173").Value = "not match"

# Row 174
$ws2.Range("A174").Value = "292-018"
$ws2.Range("B174").Value = 7
$ws2.Range("C174").Value = "In progress"
$ws2.Range("D174").Value = 7
$ws2.Range("E174").Value = 0
$ws2.Range("F174").Value = 0
$ws2.Range("G174").Value = 0
$ws2.Range("H174").Value = 7
$ws2.Range("I174").Value = "06-09-2023"
$ws2.Range("U174").Value = 1

# Row 175
$ws2.Range("A175").Value = "292-019"
$ws2.Range("B175").Value = 5
$ws2.Range("C175").Value = "In progress"
$ws2.Range("D175").Value = 5
$ws2.Range("E175").Value = 0
$ws2.Range("F175").Value = 0
$ws2.Range("G175").Value = 0
$ws2.Range("H175").Value = 5
$ws2.Range("I175").Value = "07-09-2023"
$ws2.Range("U175").Value = 1

# Row 176
$ws2.Range("A176").Value = "292-020"
$ws2.Range("B176").Value = 4
$ws2.Range("C176").Value = "In progress"
$ws2.Range("D176").Value = 4
$ws2.Range("E176").Value = 0
$ws2.Range("F176").Value = 0
$ws2.Range("G176").Value = 0
$ws2.Range("H176").Value = 4
$ws2.Range("I176").Value = "07-09-2023"
$ws2.Range("U176").Value = 1

